# Update the "Metadata" sheet: Date and Count values
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Date (B8): 2024-03-27T14:24:55+00:00 -> 2024-03-27T14:30:38+00:00
$meta.Range("B8").Value = "2024-03-27T14:30:38+00:00"

# Count (B21): 12 -> 9 ; written through TEXT() + paste-special-values so the
# result stays a genuine text cell (matching the original shared-string typed
# cell) instead of being auto-converted to a number.
$meta.Range("Z1").Formula = "=TEXT(9,""0"")"
$meta.Range("Z1").Copy()
$meta.Range("B21").PasteSpecial(-4163)
$meta.Range("Z1").ClearContents()

# Remove 3 medication rows (PROFENID, DOLIPRANE, VIT B12) from "Concepts"
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A9:D11").Delete(-4162)
